$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update predicted values in column B (rows 2-8) to the new uniform prediction
$ws.Range("B2:B8").Value = 68475.109375
